$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "added hr_periods again" - re-insert the "0.0" category row for CKD Stage.
# The existing rows 25-28 shift their category label down by one slot while
# the counts for the (re-added) "0.0" category take the values that used to
# belong to the "Absent" row.

# Keep these text-like labels stored as text (not auto-converted to numbers).
$ws.Range("B25:B28").NumberFormat = "@"

$ws.Range("B25").Value = "0.0"
$ws.Range("C25").Value = "1362 (90.6)"
$ws.Range("D25").Value = "6634 (93.8)"

$ws.Range("B26").Value = "1.0"
$ws.Range("C26").Value = "3 (0.2)"
$ws.Range("D26").Value = "1 (0.0)"

$ws.Range("B27").Value = "2.0"
$ws.Range("C27").Value = "14 (0.9)"
$ws.Range("D27").Value = "45 (0.6)"

$ws.Range("B28").Value = "3.0"
$ws.Range("C28").Value = "124 (8.3)"
$ws.Range("D28").Value = "393 (5.6)"
